# Apply transition-matrix probability updates to Sheet1
# This reflects the diff in team_specific_matrix/JWU (Charlotte)_B.xlsx
# where many zero placeholder cells were filled in with computed
# game-simulation transition probabilities (rows sum to 1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.1785714285714286
$ws.Cells.Item(2, 3).Value = 0.6071428571428571
$ws.Cells.Item(2, 16).Value = 0.1071428571428571
$ws.Cells.Item(2, 19).Value = 0.1071428571428571
$ws.Cells.Item(3, 2).Value = 0.05555555555555555
$ws.Cells.Item(3, 16).Value = 0.6111111111111112
$ws.Cells.Item(3, 19).Value = 0.3333333333333333
$ws.Cells.Item(4, 19).Value = 1
$ws.Cells.Item(6, 2).Value = 0.06666666666666667
$ws.Cells.Item(6, 6).Value = 0.06666666666666667
$ws.Cells.Item(6, 10).Value = 0.3333333333333333
$ws.Cells.Item(6, 17).Value = 0.1333333333333333
$ws.Cells.Item(6, 19).Value = 0.4
$ws.Cells.Item(7, 6).Value = 0.07692307692307693
$ws.Cells.Item(7, 10).Value = 0.3076923076923077
$ws.Cells.Item(7, 17).Value = 0.1538461538461539
$ws.Cells.Item(7, 18).Value = 0.07692307692307693
$ws.Cells.Item(7, 19).Value = 0.3846153846153846
$ws.Cells.Item(8, 2).Value = 0.18
$ws.Cells.Item(8, 4).Value = 0.02
$ws.Cells.Item(8, 6).Value = 0.06
$ws.Cells.Item(8, 10).Value = 0.12
$ws.Cells.Item(8, 15).Value = 0.04
$ws.Cells.Item(8, 17).Value = 0.2
$ws.Cells.Item(8, 18).Value = 0.02
$ws.Cells.Item(8, 19).Value = 0.36
$ws.Cells.Item(9, 2).Value = 0.1176470588235294
$ws.Cells.Item(9, 6).Value = 0.1176470588235294
$ws.Cells.Item(9, 10).Value = 0.1764705882352941
$ws.Cells.Item(9, 17).Value = 0.05882352941176471
$ws.Cells.Item(9, 19).Value = 0.5294117647058824
$ws.Cells.Item(10, 2).Value = 0.1153846153846154
$ws.Cells.Item(10, 4).Value = 0.02564102564102564
$ws.Cells.Item(10, 6).Value = 0.01282051282051282
$ws.Cells.Item(10, 10).Value = 0.1153846153846154
$ws.Cells.Item(10, 17).Value = 0.2051282051282051
$ws.Cells.Item(10, 18).Value = 0.07692307692307693
$ws.Cells.Item(10, 19).Value = 0.4487179487179487
$ws.Cells.Item(11, 7).Value = 0.125
$ws.Cells.Item(11, 10).Value = 0.0625
$ws.Cells.Item(11, 11).Value = 0.1875
$ws.Cells.Item(11, 12).Value = 0.5625
$ws.Cells.Item(11, 19).Value = 0.0625
$ws.Cells.Item(12, 7).Value = 0.6
$ws.Cells.Item(12, 10).Value = 0.2
$ws.Cells.Item(12, 11).Value = 0.1
$ws.Cells.Item(12, 12).Value = 0.1
$ws.Cells.Item(13, 7).Value = 0.8333333333333334
$ws.Cells.Item(13, 19).Value = 0.1666666666666667
$ws.Cells.Item(15, 8).Value = 0.1666666666666667
$ws.Cells.Item(15, 10).Value = 0.1666666666666667
$ws.Cells.Item(15, 19).Value = 0.6666666666666666
$ws.Cells.Item(16, 6).Value = 0.1428571428571428
$ws.Cells.Item(16, 8).Value = 0.1428571428571428
$ws.Cells.Item(16, 10).Value = 0.2857142857142857
$ws.Cells.Item(16, 11).Value = 0.1428571428571428
$ws.Cells.Item(16, 13).Value = 0.07142857142857142
$ws.Cells.Item(16, 19).Value = 0.2142857142857143
$ws.Cells.Item(17, 6).Value = 0.03333333333333333
$ws.Cells.Item(17, 8).Value = 0.2
$ws.Cells.Item(17, 9).Value = 0.1666666666666667
$ws.Cells.Item(17, 10).Value = 0.2333333333333333
$ws.Cells.Item(17, 11).Value = 0.1333333333333333
$ws.Cells.Item(17, 13).Value = 0.03333333333333333
$ws.Cells.Item(17, 14).Value = 0.03333333333333333
$ws.Cells.Item(17, 19).Value = 0.1666666666666667
$ws.Cells.Item(18, 8).Value = 0.375
$ws.Cells.Item(18, 9).Value = 0.125
$ws.Cells.Item(18, 10).Value = 0.375
$ws.Cells.Item(18, 13).Value = 0.125
$ws.Cells.Item(19, 6).Value = 0.01818181818181818
$ws.Cells.Item(19, 8).Value = 0.3454545454545455
$ws.Cells.Item(19, 9).Value = 0.1
$ws.Cells.Item(19, 10).Value = 0.3181818181818182
$ws.Cells.Item(19, 11).Value = 0.05454545454545454
$ws.Cells.Item(19, 13).Value = 0.02727272727272727
$ws.Cells.Item(19, 15).Value = 0.01818181818181818
$ws.Cells.Item(19, 19).Value = 0.1181818181818182
